$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.130287289619446
$ws.Range("B1").Value = 2.282751798629761
$ws.Range("C1").Value = 11.00678253173828
$ws.Range("D1").Value = 2.045343637466431
$ws.Range("E1").Value = 1.281050682067871
